$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SVMail")

# Update the Subject test data for the second test case row
$ws.Range("G2").Value = "Test gMail-"

# Reflect the active selection moving to G2 (as captured at save time)
$ws.Range("G2").Select()
